# NATMI ligand-receptor edge table update ("Natmi following Dr Hou advice").
#
# The Ligand-expressing-cells (col E) and Receptor-expressing-cells (col K)
# counts for every data row go from 1 -> 3 cells, which cascades through the
# already-computed NATMI specificity/weight metrics in columns G-J and M-T
# (these are literal precomputed values in the workbook, not live formulas,
# so each affected cell is written with its new literal value below).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet


# Row 2
$ws.Cells.Item(2, 5).Value = 3
$ws.Cells.Item(2, 7).Value = 18.82499266666667
$ws.Cells.Item(2, 8).Value = 56.474978
$ws.Cells.Item(2, 9).Value = 0.06886869772378311
$ws.Cells.Item(2, 10).Value = 0.0688686977237831
$ws.Cells.Item(2, 11).Value = 3
$ws.Cells.Item(2, 13).Value = 192.8285726666667
$ws.Cells.Item(2, 14).Value = 578.485718
$ws.Cells.Item(2, 15).Value = 0.7801188850698786
$ws.Cells.Item(2, 16).Value = 0.7801188850698786
$ws.Cells.Item(2, 17).Value = 3629.996466373801
$ws.Cells.Item(2, 18).Value = 32669.9681973642
$ws.Cells.Item(2, 19).Value = 0.05372577168449216
$ws.Cells.Item(2, 20).Value = 0.05372577168449216

# Row 3
$ws.Cells.Item(3, 5).Value = 3
$ws.Cells.Item(3, 7).Value = 18.82499266666667
$ws.Cells.Item(3, 8).Value = 56.474978
$ws.Cells.Item(3, 9).Value = 0.06886869772378311
$ws.Cells.Item(3, 10).Value = 0.0688686977237831
$ws.Cells.Item(3, 11).Value = 3
$ws.Cells.Item(3, 13).Value = 0.4209206666666667
$ws.Cells.Item(3, 14).Value = 1.262762
$ws.Cells.Item(3, 15).Value = 0.001702901995496819
$ws.Cells.Item(3, 16).Value = 0.001702901995496819
$ws.Cells.Item(3, 17).Value = 7.923828463248445
$ws.Cells.Item(3, 18).Value = 71.31445616923601
$ws.Cells.Item(3, 19).Value = 0.0001172766427810975
$ws.Cells.Item(3, 20).Value = 0.0001172766427810975

# Row 4
$ws.Cells.Item(4, 5).Value = 3
$ws.Cells.Item(4, 7).Value = 18.82499266666667
$ws.Cells.Item(4, 8).Value = 56.474978
$ws.Cells.Item(4, 9).Value = 0.06886869772378311
$ws.Cells.Item(4, 10).Value = 0.0688686977237831
$ws.Cells.Item(4, 11).Value = 3
$ws.Cells.Item(4, 13).Value = 45.70525533333333
$ws.Cells.Item(4, 14).Value = 137.115766
$ws.Cells.Item(4, 15).Value = 0.184907933193646
$ws.Cells.Item(4, 16).Value = 0.184907933193646
$ws.Cells.Item(4, 17).Value = 860.4010964781276
$ws.Cells.Item(4, 18).Value = 7743.609868303149
$ws.Cells.Item(4, 19).Value = 0.01273436855784269
$ws.Cells.Item(4, 20).Value = 0.01273436855784269

# Row 5
$ws.Cells.Item(5, 5).Value = 3
$ws.Cells.Item(5, 7).Value = 18.82499266666667
$ws.Cells.Item(5, 8).Value = 56.474978
$ws.Cells.Item(5, 9).Value = 0.06886869772378311
$ws.Cells.Item(5, 10).Value = 0.0688686977237831
$ws.Cells.Item(5, 11).Value = 3
$ws.Cells.Item(5, 13).Value = 8.223696
$ws.Cells.Item(5, 14).Value = 24.671088
$ws.Cells.Item(5, 15).Value = 0.0332702797409786
$ws.Cells.Item(5, 16).Value = 0.0332702797409786
$ws.Cells.Item(5, 17).Value = 154.811016892896
$ws.Cells.Item(5, 18).Value = 1393.299152036064
$ws.Cells.Item(5, 19).Value = 0.002291280838667161
$ws.Cells.Item(5, 20).Value = 0.00229128083866716

# Row 6
$ws.Cells.Item(6, 5).Value = 3
$ws.Cells.Item(6, 7).Value = 121.8208923333333
$ws.Cells.Item(6, 8).Value = 365.462677
$ws.Cells.Item(6, 9).Value = 0.4456653109566078
$ws.Cells.Item(6, 10).Value = 0.4456653109566078
$ws.Cells.Item(6, 11).Value = 3
$ws.Cells.Item(6, 13).Value = 192.8285726666667
$ws.Cells.Item(6, 14).Value = 578.485718
$ws.Cells.Item(6, 15).Value = 0.7801188850698786
$ws.Cells.Item(6, 16).Value = 0.7801188850698786
$ws.Cells.Item(6, 17).Value = 23490.54878961635
$ws.Cells.Item(6, 18).Value = 211414.9391065471
$ws.Cells.Item(6, 19).Value = 0.3476719254977896
$ws.Cells.Item(6, 20).Value = 0.3476719254977896

# Row 7
$ws.Cells.Item(7, 5).Value = 3
$ws.Cells.Item(7, 7).Value = 121.8208923333333
$ws.Cells.Item(7, 8).Value = 365.462677
$ws.Cells.Item(7, 9).Value = 0.4456653109566078
$ws.Cells.Item(7, 10).Value = 0.4456653109566078
$ws.Cells.Item(7, 11).Value = 3
$ws.Cells.Item(7, 13).Value = 0.4209206666666667
$ws.Cells.Item(7, 14).Value = 1.262762
$ws.Cells.Item(7, 15).Value = 0.001702901995496819
$ws.Cells.Item(7, 16).Value = 0.001702901995496819
$ws.Cells.Item(7, 17).Value = 51.2769312148749
$ws.Cells.Item(7, 18).Value = 461.4923809338741
$ws.Cells.Item(7, 19).Value = 0.0007589243473517178
$ws.Cells.Item(7, 20).Value = 0.0007589243473517176

# Row 8
$ws.Cells.Item(8, 5).Value = 3
$ws.Cells.Item(8, 7).Value = 121.8208923333333
$ws.Cells.Item(8, 8).Value = 365.462677
$ws.Cells.Item(8, 9).Value = 0.4456653109566078
$ws.Cells.Item(8, 10).Value = 0.4456653109566078
$ws.Cells.Item(8, 11).Value = 3
$ws.Cells.Item(8, 13).Value = 45.70525533333333
$ws.Cells.Item(8, 14).Value = 137.115766
$ws.Cells.Item(8, 15).Value = 0.184907933193646
$ws.Cells.Item(8, 16).Value = 0.184907933193646
$ws.Cells.Item(8, 17).Value = 5567.85498902951
$ws.Cells.Item(8, 18).Value = 50110.69490126559
$ws.Cells.Item(8, 19).Value = 0.08240705154508993
$ws.Cells.Item(8, 20).Value = 0.08240705154508991

# Row 9
$ws.Cells.Item(9, 5).Value = 3
$ws.Cells.Item(9, 7).Value = 121.8208923333333
$ws.Cells.Item(9, 8).Value = 365.462677
$ws.Cells.Item(9, 9).Value = 0.4456653109566078
$ws.Cells.Item(9, 10).Value = 0.4456653109566078
$ws.Cells.Item(9, 11).Value = 3
$ws.Cells.Item(9, 13).Value = 8.223696
$ws.Cells.Item(9, 14).Value = 24.671088
$ws.Cells.Item(9, 15).Value = 0.0332702797409786
$ws.Cells.Item(9, 16).Value = 0.0332702797409786
$ws.Cells.Item(9, 17).Value = 1001.817984998064
$ws.Cells.Item(9, 18).Value = 9016.361864982577
$ws.Cells.Item(9, 19).Value = 0.01482740956637656
$ws.Cells.Item(9, 20).Value = 0.01482740956637655

# Row 10
$ws.Cells.Item(10, 5).Value = 3
$ws.Cells.Item(10, 7).Value = 87.673585
$ws.Cells.Item(10, 8).Value = 263.020755
$ws.Cells.Item(10, 9).Value = 0.3207419907481189
$ws.Cells.Item(10, 10).Value = 0.3207419907481188
$ws.Cells.Item(10, 11).Value = 3
$ws.Cells.Item(10, 13).Value = 192.8285726666667
$ws.Cells.Item(10, 14).Value = 578.485718
$ws.Cells.Item(10, 15).Value = 0.7801188850698786
$ws.Cells.Item(10, 16).Value = 0.7801188850698786
$ws.Cells.Item(10, 17).Value = 16905.97225611968
$ws.Cells.Item(10, 18).Value = 152153.7503050771
$ws.Cells.Item(10, 19).Value = 0.2502168842175158
$ws.Cells.Item(10, 20).Value = 0.2502168842175158

# Row 11
$ws.Cells.Item(11, 5).Value = 3
$ws.Cells.Item(11, 7).Value = 87.673585
$ws.Cells.Item(11, 8).Value = 263.020755
$ws.Cells.Item(11, 9).Value = 0.3207419907481189
$ws.Cells.Item(11, 10).Value = 0.3207419907481188
$ws.Cells.Item(11, 11).Value = 3
$ws.Cells.Item(11, 13).Value = 0.4209206666666667
$ws.Cells.Item(11, 14).Value = 1.262762
$ws.Cells.Item(11, 15).Value = 0.001702901995496819
$ws.Cells.Item(11, 16).Value = 0.001702901995496819
$ws.Cells.Item(11, 17).Value = 36.90362384725668
$ws.Cells.Item(11, 18).Value = 332.1326146253101
$ws.Cells.Item(11, 19).Value = 0.0005461921760845939
$ws.Cells.Item(11, 20).Value = 0.0005461921760845938

# Row 12
$ws.Cells.Item(12, 5).Value = 3
$ws.Cells.Item(12, 7).Value = 87.673585
$ws.Cells.Item(12, 8).Value = 263.020755
$ws.Cells.Item(12, 9).Value = 0.3207419907481189
$ws.Cells.Item(12, 10).Value = 0.3207419907481188
$ws.Cells.Item(12, 11).Value = 3
$ws.Cells.Item(12, 13).Value = 45.70525533333333
$ws.Cells.Item(12, 14).Value = 137.115766
$ws.Cells.Item(12, 15).Value = 0.184907933193646
$ws.Cells.Item(12, 16).Value = 0.184907933193646
$ws.Cells.Item(12, 17).Value = 4007.143588413704
$ws.Cells.Item(12, 18).Value = 36064.29229572333
$ws.Cells.Item(12, 19).Value = 0.0593077385976502
$ws.Cells.Item(12, 20).Value = 0.05930773859765019

# Row 13
$ws.Cells.Item(13, 5).Value = 3
$ws.Cells.Item(13, 7).Value = 87.673585
$ws.Cells.Item(13, 8).Value = 263.020755
$ws.Cells.Item(13, 9).Value = 0.3207419907481189
$ws.Cells.Item(13, 10).Value = 0.3207419907481188
$ws.Cells.Item(13, 11).Value = 3
$ws.Cells.Item(13, 13).Value = 8.223696
$ws.Cells.Item(13, 14).Value = 24.671088
$ws.Cells.Item(13, 15).Value = 0.0332702797409786
$ws.Cells.Item(13, 16).Value = 0.0332702797409786
$ws.Cells.Item(13, 17).Value = 721.00091027016
$ws.Cells.Item(13, 18).Value = 6489.008192431441
$ws.Cells.Item(13, 19).Value = 0.01067117575686829
$ws.Cells.Item(13, 20).Value = 0.01067117575686828

# Row 14
$ws.Cells.Item(14, 5).Value = 3
$ws.Cells.Item(14, 7).Value = 45.02666966666666
$ws.Cells.Item(14, 8).Value = 135.080009
$ws.Cells.Item(14, 9).Value = 0.1647240005714903
$ws.Cells.Item(14, 10).Value = 0.1647240005714903
$ws.Cells.Item(14, 11).Value = 3
$ws.Cells.Item(14, 13).Value = 192.8285726666667
$ws.Cells.Item(14, 14).Value = 578.485718
$ws.Cells.Item(14, 15).Value = 0.7801188850698786
$ws.Cells.Item(14, 16).Value = 0.7801188850698786
$ws.Cells.Item(14, 17).Value = 8682.428443756829
$ws.Cells.Item(14, 18).Value = 78141.85599381146
$ws.Cells.Item(14, 19).Value = 0.128504303670081
$ws.Cells.Item(14, 20).Value = 0.128504303670081

# Row 15
$ws.Cells.Item(15, 5).Value = 3
$ws.Cells.Item(15, 7).Value = 45.02666966666666
$ws.Cells.Item(15, 8).Value = 135.080009
$ws.Cells.Item(15, 9).Value = 0.1647240005714903
$ws.Cells.Item(15, 10).Value = 0.1647240005714903
$ws.Cells.Item(15, 11).Value = 3
$ws.Cells.Item(15, 13).Value = 0.4209206666666667
$ws.Cells.Item(15, 14).Value = 1.262762
$ws.Cells.Item(15, 15).Value = 0.001702901995496819
$ws.Cells.Item(15, 16).Value = 0.001702901995496819
$ws.Cells.Item(15, 17).Value = 18.95265581387311
$ws.Cells.Item(15, 18).Value = 170.573902324858
$ws.Cells.Item(15, 19).Value = 0.00028050882927941
$ws.Cells.Item(15, 20).Value = 0.0002805088292794098

# Row 16
$ws.Cells.Item(16, 5).Value = 3
$ws.Cells.Item(16, 7).Value = 45.02666966666666
$ws.Cells.Item(16, 8).Value = 135.080009
$ws.Cells.Item(16, 9).Value = 0.1647240005714903
$ws.Cells.Item(16, 10).Value = 0.1647240005714903
$ws.Cells.Item(16, 11).Value = 3
$ws.Cells.Item(16, 13).Value = 45.70525533333333
$ws.Cells.Item(16, 14).Value = 137.115766
$ws.Cells.Item(16, 15).Value = 0.184907933193646
$ws.Cells.Item(16, 16).Value = 0.184907933193646
$ws.Cells.Item(16, 17).Value = 2057.955433924655
$ws.Cells.Item(16, 18).Value = 18521.59890532189
$ws.Cells.Item(16, 19).Value = 0.03045877449306324
$ws.Cells.Item(16, 20).Value = 0.03045877449306323

# Row 17
$ws.Cells.Item(17, 5).Value = 3
$ws.Cells.Item(17, 7).Value = 45.02666966666666
$ws.Cells.Item(17, 8).Value = 135.080009
$ws.Cells.Item(17, 9).Value = 0.1647240005714903
$ws.Cells.Item(17, 10).Value = 0.1647240005714903
$ws.Cells.Item(17, 11).Value = 3
$ws.Cells.Item(17, 13).Value = 8.223696
$ws.Cells.Item(17, 14).Value = 24.671088
$ws.Cells.Item(17, 15).Value = 0.0332702797409786
$ws.Cells.Item(17, 16).Value = 0.0332702797409786
$ws.Cells.Item(17, 17).Value = 370.285643231088
$ws.Cells.Item(17, 18).Value = 3332.570789079792
$ws.Cells.Item(17, 19).Value = 0.005480413579066601
$ws.Cells.Item(17, 20).Value = 0.005480413579066599

Write-Host "Updated ligand/receptor-expressing cell counts (1 -> 3) and recomputed NATMI metrics for rows 2-17."